# "Generate Report for Handoff"
#
# The workbook tracks localization hand-off/hand-back status for two
# source files (883008bc... and e8172bf7...) across three sheets:
#   Overview (summary), zh-cn (detail), de-de (detail).
#
# This edit:
#   1. Flips the 883008bc... entry's status from
#      "Handed back: in sync with en-US" to "Ready for handoff" (and bumps
#      its "Latest Handoff Datetime" on the zh-cn/de-de detail sheets,
#      since a fresh hand-off was just generated).
#   2. Drops the e8172bf7... entry entirely (its row) from all three sheets
#      - it shipped/handed-back earlier and is no longer tracked here.
#
# Row deletion in this engine does not renumber/prune the <hyperlinks>
# collection, so hyperlinks are rebuilt from scratch per sheet afterwards.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"

# Remove the e8172bf7 row (row 3); .localization-config shifts up to row 3.
$ov.Rows.Item(3).Delete()

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8c2d4629df2344af058dcbef92838bb6665d30b7/e2e/883008bc-1e48-4fbd-aa59-c5b93c4a3497.md", "", "", "883008bc-1e48-4fbd-aa59-c5b93c4a3497.md")
$ov.Hyperlinks.Add($ov.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8c2d4629df2344af058dcbef92838bb6665d30b7/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("B2").Value = "Ready for handoff"
$zh.Range("D2").Value = "2016-03-10 09:26:20"

# Remove the e8172bf7 row (row 3); .localization-config shifts up to row 3.
$zh.Rows.Item(3).Delete()

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8c2d4629df2344af058dcbef92838bb6665d30b7/e2e/883008bc-1e48-4fbd-aa59-c5b93c4a3497.md", "", "", "883008bc-1e48-4fbd-aa59-c5b93c4a3497.md")
$zh.Hyperlinks.Add($zh.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/dc68aae20cc732b53f8f06124d556c2d304c4909/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/883008bc-1e48-4fbd-aa59-c5b93c4a3497.d996985f1bc82340bc7808d93a440a7a7776aaaf.zh-cn.xlf", "", "", "883008bc-1e48-4fbd-aa59-c5b93c4a3497.d996985f1bc82340bc7808d93a440a7a7776aaaf.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/5793675f66f47f7b0605797f4a8f9145cf7e9dc5/e2e/883008bc-1e48-4fbd-aa59-c5b93c4a3497.md", "", "", "883008bc-1e48-4fbd-aa59-c5b93c4a3497.md")
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c4b391983839071a68e182ee1453f4d7ef333077/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/883008bc-1e48-4fbd-aa59-c5b93c4a3497.d996985f1bc82340bc7808d93a440a7a7776aaaf.zh-cn.xlf", "", "", "883008bc-1e48-4fbd-aa59-c5b93c4a3497.d996985f1bc82340bc7808d93a440a7a7776aaaf.zh-cn.xlf")
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8c2d4629df2344af058dcbef92838bb6665d30b7/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("B2").Value = "Ready for handoff"
$de.Range("D2").Value = "2016-03-10 09:26:28"

# Remove the e8172bf7 row (row 3); .localization-config shifts up to row 3.
$de.Rows.Item(3).Delete()

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/8c2d4629df2344af058dcbef92838bb6665d30b7/e2e/883008bc-1e48-4fbd-aa59-c5b93c4a3497.md", "", "", "883008bc-1e48-4fbd-aa59-c5b93c4a3497.md")
$de.Hyperlinks.Add($de.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/158e2c403370ad5857c7eb44d5254916abd1d0e5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/883008bc-1e48-4fbd-aa59-c5b93c4a3497.d996985f1bc82340bc7808d93a440a7a7776aaaf.de-de.xlf", "", "", "883008bc-1e48-4fbd-aa59-c5b93c4a3497.d996985f1bc82340bc7808d93a440a7a7776aaaf.de-de.xlf")
$de.Hyperlinks.Add($de.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/98a6bb084466de1e3c78866a5bc20315741a3a5d/e2e/883008bc-1e48-4fbd-aa59-c5b93c4a3497.md", "", "", "883008bc-1e48-4fbd-aa59-c5b93c4a3497.md")
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/73a0437e32d2825a3ccec0df9d37b3b5fa41aae7/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/883008bc-1e48-4fbd-aa59-c5b93c4a3497.d996985f1bc82340bc7808d93a440a7a7776aaaf.de-de.xlf", "", "", "883008bc-1e48-4fbd-aa59-c5b93c4a3497.d996985f1bc82340bc7808d93a440a7a7776aaaf.de-de.xlf")
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8c2d4629df2344af058dcbef92838bb6665d30b7/.localization-config", "", "", ".localization-config")
